$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1041.2
$ws.Range("I19").Value = 665
$ws.Range("J19").Value = 1292
$ws.Range("K19").Value = 665
$ws.Range("L19").Value = 1292
$ws.Range("M19").Value = -490
$ws.Range("N19").Value = -1642
$ws.Range("H116").Value = 3253.75
$ws.Range("I116").Value = 4185.5713
$ws.Range("J116").Value = 1949.2
$ws.Range("K116").Value = 4185.5713
$ws.Range("L116").Value = 1949.2
$ws.Range("M116").Value = -743.5712999999996
$ws.Range("N116").Value = -8833.200000000001
$ws.Range("H138").Value = 6216409
$ws.Range("I138").Value = 7148246
$ws.Range("J138").Value = 6102770.5
$ws.Range("K138").Value = 21444738
$ws.Range("L138").Value = 18308311.5
$ws.Range("M138").Value = -21439598
$ws.Range("N138").Value = -18318591.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1708.8462
$ws.Range("I45").Value = 1733.4166
$ws.Range("J45").Value = 1414
$ws.Range("K45").Value = 1733.4166
$ws.Range("L45").Value = 1414
$ws.Range("M45").Value = -1356.4166
$ws.Range("N45").Value = -2168
$ws.Range("H74").Value = 5337.3125
$ws.Range("I74").Value = 1053.6364
$ws.Range("J74").Value = 14761.4
$ws.Range("K74").Value = 1053.6364
$ws.Range("L74").Value = 14761.4
$ws.Range("M74").Value = -179.6364000000001
$ws.Range("N74").Value = -16509.4
$ws.Range("H77").Value = 5337.3125
$ws.Range("I77").Value = 1053.6364
$ws.Range("J77").Value = 14761.4
$ws.Range("K77").Value = 5268.182000000001
$ws.Range("L77").Value = 73807
$ws.Range("M77").Value = -900.1820000000007
$ws.Range("N77").Value = -82543
$ws.Range("H122").Value = 1717.5834
$ws.Range("I122").Value = 1567.8889
$ws.Range("K122").Value = 4703.6667
$ws.Range("M122").Value = -2253.6667
$ws.Range("H133").Value = 42632.375
$ws.Range("J133").Value = 42632.375
$ws.Range("L133").Value = 42632.375
$ws.Range("N133").Value = -47692.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1899.6364
$ws.Range("I86").Value = 2347.25
$ws.Range("J86").Value = 1643.8572
$ws.Range("K86").Value = 2347.25
$ws.Range("L86").Value = 1643.8572
$ws.Range("M86").Value = -1224.25
$ws.Range("N86").Value = -3889.8572
$ws.Range("H89").Value = 1899.6364
$ws.Range("I89").Value = 2347.25
$ws.Range("J89").Value = 1643.8572
$ws.Range("K89").Value = 11736.25
$ws.Range("L89").Value = 8219.286
$ws.Range("M89").Value = -6120.25
$ws.Range("N89").Value = -19451.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2001.4286
$ws.Range("I16").Value = 1902
$ws.Range("J16").Value = 2250
$ws.Range("K16").Value = 1902
$ws.Range("L16").Value = 2250
$ws.Range("M16").Value = -1615
$ws.Range("N16").Value = -2824
$ws.Range("H31").Value = 5901.6787
$ws.Range("I31").Value = 2775.4167
$ws.Range("J31").Value = 8246.375
$ws.Range("K31").Value = 2775.4167
$ws.Range("L31").Value = 8246.375
$ws.Range("M31").Value = -2480.4167
$ws.Range("N31").Value = -8836.375
$ws.Range("H34").Value = 5901.6787
$ws.Range("I34").Value = 2775.4167
$ws.Range("J34").Value = 8246.375
$ws.Range("K34").Value = 2775.4167
$ws.Range("L34").Value = 8246.375
$ws.Range("M34").Value = -2573.4167
$ws.Range("N34").Value = -8650.375
$ws.Range("H58").Value = 25642728
$ws.Range("I58").Value = 35715584
$ws.Range("J58").Value = 2730.3635
$ws.Range("K58").Value = 35715584
$ws.Range("L58").Value = 2730.3635
$ws.Range("M58").Value = -35715381
$ws.Range("N58").Value = -3136.3635
$ws.Range("H102").Value = 29900
$ws.Range("J102").Value = 29900
$ws.Range("L102").Value = 29900
$ws.Range("N102").Value = -34768
$ws.Range("H113").Value = 2001.4286
$ws.Range("I113").Value = 1902
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 1902
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 268
$ws.Range("N113").Value = -6590
$ws.Range("H136").Value = 25642728
$ws.Range("I136").Value = 35715584
$ws.Range("J136").Value = 2730.3635
$ws.Range("K136").Value = 107146752
$ws.Range("L136").Value = 8191.0905
$ws.Range("M136").Value = -107144202
$ws.Range("N136").Value = -13291.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3464.6155
$ws.Range("I3").Value = 1137.7778
$ws.Range("J3").Value = 8700
$ws.Range("K3").Value = 3413.3334
$ws.Range("L3").Value = 26100
$ws.Range("M3").Value = -3301.3334
$ws.Range("N3").Value = -26324
$ws.Range("H5").Value = 2336.4
$ws.Range("I5").Value = 946.1539
$ws.Range("J5").Value = 3842.5
$ws.Range("K5").Value = 2838.4617
$ws.Range("L5").Value = 11527.5
$ws.Range("M5").Value = -2726.4617
$ws.Range("N5").Value = -11751.5
$ws.Range("H68").Value = 280500.5
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 280500.5
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112
$ws.Range("H76").Value = 2950
$ws.Range("I76").Value = 2950
$ws.Range("K76").Value = 8850
$ws.Range("M76").Value = -8467
$ws.Range("H79").Value = 2950
$ws.Range("I79").Value = 2950
$ws.Range("K79").Value = 8850
$ws.Range("M79").Value = -7524
$ws.Range("H97").Value = 1083.8462
$ws.Range("I97").Value = 999
$ws.Range("J97").Value = 1182.8334
$ws.Range("K97").Value = 2997
$ws.Range("L97").Value = 3548.5002
$ws.Range("M97").Value = -2501
$ws.Range("N97").Value = -4540.5002
$ws.Range("H112").Value = 111113370
$ws.Range("I112").Value = 1133.1666
$ws.Range("J112").Value = 333337820
$ws.Range("K112").Value = 3399.4998
$ws.Range("L112").Value = 1000013460
$ws.Range("M112").Value = -2291.4998
$ws.Range("N112").Value = -1000015676
$ws.Range("H122").Value = 1283159.5
$ws.Range("I122").Value = 679.8
$ws.Range("J122").Value = 1710652.8
$ws.Range("K122").Value = 6118.2
$ws.Range("L122").Value = 15395875.2
$ws.Range("M122").Value = -3668.2
$ws.Range("N122").Value = -15400775.2
$ws.Range("H129").Value = 1302.1818
$ws.Range("J129").Value = 1918.9231
$ws.Range("L129").Value = 5756.7693
$ws.Range("N129").Value = -15756.7693
$ws.Range("H131").Value = 4445874.5
$ws.Range("I131").Value = 756.6667
$ws.Range("J131").Value = 4631088
$ws.Range("K131").Value = 2270.0001
$ws.Range("L131").Value = 13893264
$ws.Range("M131").Value = 2769.9999
$ws.Range("N131").Value = -13903344
$ws.Range("H134").Value = 6915.2
$ws.Range("I134").Value = 4272
$ws.Range("J134").Value = 10880
$ws.Range("K134").Value = 12816
$ws.Range("L134").Value = 32640
$ws.Range("M134").Value = -7746
$ws.Range("N134").Value = -42780
$ws.Range("H135").Value = 2336.4
$ws.Range("I135").Value = 946.1539
$ws.Range("J135").Value = 3842.5
$ws.Range("K135").Value = 8515.3851
$ws.Range("L135").Value = 34582.5
$ws.Range("M135").Value = -5980.3851
$ws.Range("N135").Value = -39652.5
$ws.Range("H137").Value = 4952
$ws.Range("I137").Value = 5377.5
$ws.Range("J137").Value = 3250
$ws.Range("K137").Value = 16132.5
$ws.Range("L137").Value = 9750
$ws.Range("M137").Value = -11032.5
$ws.Range("N137").Value = -19950
$ws.Range("H139").Value = 4369.9
$ws.Range("I139").Value = 4216.125
$ws.Range("J139").Value = 4985
$ws.Range("K139").Value = 12648.375
$ws.Range("L139").Value = 14955
$ws.Range("M139").Value = -7508.375
$ws.Range("N139").Value = -25235
$ws.Range("H140").Value = 7319.775
$ws.Range("I140").Value = 11519
$ws.Range("J140").Value = 3520.476
$ws.Range("K140").Value = 34557
$ws.Range("L140").Value = 10561.428
$ws.Range("M140").Value = -29377
$ws.Range("N140").Value = -20921.428
$ws.Range("H141").Value = 6033.75
$ws.Range("I141").Value = 8354
$ws.Range("J141").Value = 2166.6667
$ws.Range("K141").Value = 25062
$ws.Range("L141").Value = 6500.000100000001
$ws.Range("M141").Value = -19882
$ws.Range("N141").Value = -16860.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 18310.566
$ws.Range("J121").Value = 18310.566
$ws.Range("L121").Value = 18310.566
$ws.Range("N121").Value = -21804.566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3842.1428
$ws.Range("I122").Value = 3790
$ws.Range("J122").Value = 3846.1538
$ws.Range("K122").Value = 11370
$ws.Range("L122").Value = 11538.4614
$ws.Range("M122").Value = -8920
$ws.Range("N122").Value = -16438.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 199666.67
$ws.Range("J46").Value = 199666.67
$ws.Range("L46").Value = 199666.67
$ws.Range("N46").Value = -200128.67
$ws.Range("H126").Value = 587.86365
$ws.Range("I126").Value = 401.8421
$ws.Range("J126").Value = 1766
$ws.Range("K126").Value = 1205.5263
$ws.Range("L126").Value = 5298
$ws.Range("M126").Value = 1264.4737
$ws.Range("N126").Value = -10238
$ws.Range("H134").Value = 199666.67
$ws.Range("J134").Value = 199666.67
$ws.Range("L134").Value = 599000.01
$ws.Range("N134").Value = -604070.01
